$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Metadata")
$ws2 = $wb.Worksheets.Item("Concepts")

# Update the Date value (row 8, column B on the Metadata sheet)
$ws1.Range("B8").Value = "2024-03-15T21:50:53+00:00"

# Update the Count value (row 21, column B on the Metadata sheet) - force
# text storage (to match the source data's shared-string "10") while
# keeping the original cell formatting.
$ws1.Range("B21").Value = "'10"
$ws1.Range("B20").Copy()
$ws1.Range("B21").PasteSpecial(-4122)

# Add the new "GSR" concept row (row 11) on the Concepts sheet, matching
# the existing row formatting/types exactly.
$ws2.Range("A10").Copy()
$ws2.Range("A11").PasteSpecial(-4104)
$ws2.Range("A10:D10").Copy()
$ws2.Range("A11:D11").PasteSpecial(-4122)
$ws2.Range("B11").Value = "GSR"
$ws2.Range("C11").Value = "Genomic Summary Results"
